$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.133099999999998
$ws.Range("A10").Value = -21.8939
$ws.Range("A12").Value = -21.53790000000001
$ws.Range("D15").Value = -8.044599999999994
$ws.Range("A18").Value = -22.44110000000002
$ws.Range("D20").Value = -7.824199999999998
$ws.Range("D29").Value = -7.235400000000001
$ws.Range("D30").Value = -7.360700000000005
$ws.Range("D31").Value = -8.439899999999993
$ws.Range("A37").Value = -19.5807
$ws.Range("D40").Value = -8.119899999999994
$ws.Range("A55").Value = -22.22980000000001
$ws.Range("A68").Value = -21.5371
$ws.Range("D68").Value = -6.941199999999998
$ws.Range("D76").Value = -7.247399999999998
$ws.Range("A77").Value = -20.92049999999999
$ws.Range("A78").Value = -20.57539999999999
$ws.Range("D87").Value = -7.892399999999993
$ws.Range("D88").Value = -7.085599999999994
$ws.Range("D96").Value = -7.528500000000004
$ws.Range("D98").Value = -8.3299
$ws.Range("D101").Value = -7.815100000000002
$ws.Range("D102").Value = -7.8497
